# Saldo.xlsx update script
# Applies the row-level changes described by the upstream diff:
#  - Update DANIELA's (004001621) balance
#  - Remove FRANCISCO (004567324) and AYRTON (001000882) rows
#  - Insert ELSI (004855073) and VALERIA (004328934 / 18429.12) rows before ELISA
#  - Insert MATEUS (004451652) row before NADY
#  - Insert PRISCILLA (004224284) row before ALESSANDRA
#  - Insert GUILHERME (004574428 / 2942.61) row before DAIANNE
#  - Remove CAIO (004512434) row
#  - Insert TIAGO (004498637) row before LEONE
#  - Remove the duplicate GUILHERME (004574428 / 596.31) row
#  - Remove the duplicate VALERIA (004328934 / 583.29) row
#
# Operations are applied from the bottom of the sheet upward so that row
# numbers used below remain valid (inserting/deleting rows only shifts the
# rows that come after them).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force the value to be stored as text so leading zeros (account
    # numbers) survive the round trip, mirroring how the other account
    # numbers in the sheet are already stored as inline/shared strings.
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

function Insert-DataRow($row, $conta, $nome, $saldo) {
    $ws.Rows.Item($row).Insert()
    Set-TextCell $row 1 $conta
    $ws.Cells.Item($row, 2).Value = $nome
    $ws.Cells.Item($row, 3).Value = $saldo
}

# --- Bottom-most edits first ---

# Remove VALERIA 004328934 / 583.29 (row 47, right after CARLA 004806244)
$ws.Rows.Item(47).Delete()

# Remove GUILHERME 004574428 / 596.31 (row 45, right before CARLA 004806244)
$ws.Rows.Item(45).Delete()

# Insert TIAGO 004498637 / 1000 before LEONE (row 22)
Insert-DataRow 22 "004498637" "TIAGO" 1000

# Remove CAIO 004512434 / 1510 (row 18)
$ws.Rows.Item(18).Delete()

# Insert GUILHERME 004574428 / 2942.61 before DAIANNE (row 16)
Insert-DataRow 16 "004574428" "GUILHERME" 2942.61

# Insert PRISCILLA 004224284 / 10000 before ALESSANDRA (row 9)
Insert-DataRow 9 "004224284" "PRISCILLA" 10000

# Insert MATEUS 004451652 / 13279.61 before NADY (row 8)
Insert-DataRow 8 "004451652" "MATEUS" 13279.61

# Insert ELSI 004855073 / 20641.08 and VALERIA 004328934 / 18429.12 before ELISA (row 6)
Insert-DataRow 6 "004855073" "ELSI" 20641.08
Insert-DataRow 7 "004328934" "VALERIA" 18429.12

# Remove FRANCISCO 004567324 (row 3) and AYRTON 001000882 (row 4, becomes row 3 after first delete)
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()

# Update DANIELA's (004001621) balance on row 2
$ws.Cells.Item(2, 3).Value = 102010.67

Write-Host "Saldo.xlsx updates applied"
